$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.39
$ws.Range("AC2").Value = 10
$ws.Range("I3").Value = 3.95
$ws.Range("N3").Value = 3.45
$ws.Range("V3").Value = 1.37
$ws.Range("W3").Value = 1.64
$ws.Range("G4").Value = 1.46
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 2.18
$ws.Range("X4").Value = 38
$ws.Range("Y4").Value = 46
$ws.Range("AB4").Value = 14.5
$ws.Range("AC4").Value = 16.5
$ws.Range("AD4").Value = 40
$ws.Range("AF4").Value = 12.5
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 29
$ws.Range("AJ4").Value = 15
$ws.Range("AK4").Value = 17
$ws.Range("AL4").Value = 36
$ws.Range("AN4").Value = 5.4
$ws.Range("AC5").Value = 970
$ws.Range("N6").Value = 3.2
$ws.Range("P6").Value = 1.74
$ws.Range("W7").Value = 2.82
$ws.Range("AB7").Value = 970
$ws.Range("AC7").Value = 970
$ws.Range("AG7").Value = 970
$ws.Range("AJ7").Value = 970
$ws.Range("AK7").Value = 970
$ws.Range("K8").Value = 6.4
$ws.Range("P8").Value = 2.22
$ws.Range("R8").Value = 1.51
$ws.Range("S8").Value = 2.56
$ws.Range("Y10").Value = 970
$ws.Range("AA10").Value = 970
$ws.Range("AC10").Value = 970
$ws.Range("AE10").Value = 970
$ws.Range("F11").Value = 1.54
$ws.Range("G11").Value = 1.56
$ws.Range("H11").Value = 7.8
$ws.Range("I11").Value = 9.4
$ws.Range("K11").Value = 4.4
$ws.Range("N11").Value = 2.98
$ws.Range("P11").Value = 1.68
$ws.Range("V11").Value = 1.11
$ws.Range("AH11").Value = 36
$ws.Range("AI11").Value = 210
$ws.Range("AM11").Value = 290
$ws.Range("L13").Value = 1.51
$ws.Range("I15").Value = 2.34
$ws.Range("L15").Value = 1.43
$ws.Range("M15").Value = 1.1
$ws.Range("S15").Value = 4.6
$ws.Range("T15").Value = 2
$ws.Range("U15").Value = 1.81
$ws.Range("V15").Value = 1.77
$ws.Range("AB16").Value = 970
$ws.Range("AC16").Value = 970
$ws.Range("AF16").Value = 970
$ws.Range("AG16").Value = 970
$ws.Range("AJ16").Value = 970
$ws.Range("G17").Value = 1.52
$ws.Range("W17").Value = 2.92
$ws.Range("L18").Value = 1.41
$ws.Range("N18").Value = 2.96
$ws.Range("O18").Value = 1.42
$ws.Range("P18").Value = 1.66
$ws.Range("Q18").Value = 2.04
$ws.Range("R18").Value = 1.24
$ws.Range("S18").Value = 4.2
$ws.Range("T18").Value = 1.8
$ws.Range("U18").Value = 1.73
$ws.Range("F19").Value = 2.22
$ws.Range("F20").Value = 1.9
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 4.1
$ws.Range("I20").Value = 4.6
$ws.Range("J20").Value = 3.75
$ws.Range("K20").Value = 3.95
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 1.73
$ws.Range("V20").Value = 1.28
$ws.Range("W20").Value = 2
$ws.Range("AD20").Value = 22
$ws.Range("AF20").Value = 13.5
$ws.Range("AH20").Value = 22
$ws.Range("F21").Value = 1.86
$ws.Range("T21").Value = 1.61
$ws.Range("U21").Value = 2.56
$ws.Range("AG21").Value = 9.800000000000001
$ws.Range("L22").Value = 1.42
$ws.Range("H23").Value = 16.5
$ws.Range("K23").Value = 8.4
$ws.Range("P23").Value = 2.34
$ws.Range("R23").Value = 1.56
$ws.Range("AL23").Value = 55
